$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1547169811320755
$ws.Range("C2").Value = 0.6339622641509434
$ws.Range("J2").Value = 0.007547169811320755
$ws.Range("P2").Value = 0.1056603773584906
$ws.Range("S2").Value = 0.09811320754716982
$ws.Range("B3").Value = 0.01176470588235294
$ws.Range("C3").Value = 0.005882352941176471
$ws.Range("J3").Value = 0.01176470588235294
$ws.Range("P3").Value = 0.7470588235294118
$ws.Range("S3").Value = 0.2235294117647059
$ws.Range("J4").Value = 0.01639344262295082
$ws.Range("P4").Value = 0.6885245901639344
$ws.Range("S4").Value = 0.2950819672131147
$ws.Range("B6").Value = 0.04950495049504951
$ws.Range("D6").Value = 0.01485148514851485
$ws.Range("E6").Value = 0.004950495049504951
$ws.Range("F6").Value = 0.09405940594059406
$ws.Range("J6").Value = 0.2227722772277228
$ws.Range("O6").Value = 0.0198019801980198
$ws.Range("Q6").Value = 0.1485148514851485
$ws.Range("R6").Value = 0.06930693069306931
$ws.Range("S6").Value = 0.3762376237623762
$ws.Range("B7").Value = 0.08496732026143791
$ws.Range("D7").Value = 0.0392156862745098
$ws.Range("F7").Value = 0.0130718954248366
$ws.Range("J7").Value = 0.1764705882352941
$ws.Range("O7").Value = 0.0196078431372549
$ws.Range("Q7").Value = 0.1568627450980392
$ws.Range("R7").Value = 0.0915032679738562
$ws.Range("S7").Value = 0.4183006535947713
$ws.Range("B8").Value = 0.1063829787234043
$ws.Range("D8").Value = 0.03073286052009456
$ws.Range("F8").Value = 0.0425531914893617
$ws.Range("J8").Value = 0.1040189125295508
$ws.Range("O8").Value = 0.01891252955082742
$ws.Range("Q8").Value = 0.2080378250591016
$ws.Range("R8").Value = 0.09456264775413711
$ws.Range("S8").Value = 0.3947990543735225
$ws.Range("B9").Value = 0.05154639175257732
$ws.Range("D9").Value = 0.03608247422680412
$ws.Range("F9").Value = 0.04639175257731959
$ws.Range("J9").Value = 0.07216494845360824
$ws.Range("O9").Value = 0.005154639175257732
$ws.Range("Q9").Value = 0.2061855670103093
$ws.Range("R9").Value = 0.1288659793814433
$ws.Range("S9").Value = 0.4536082474226804
$ws.Range("B10").Value = 0.116410670978173
$ws.Range("D10").Value = 0.02667744543249798
$ws.Range("E10").Value = 0.001616814874696847
$ws.Range("F10").Value = 0.068714632174616
$ws.Range("J10").Value = 0.09943411479385611
$ws.Range("O10").Value = 0.01616814874696847
$ws.Range("Q10").Value = 0.2021018593371059
$ws.Range("R10").Value = 0.06790622473726758
$ws.Range("S10").Value = 0.4009700889248181
$ws.Range("G11").Value = 0.125
$ws.Range("J11").Value = 0.08333333333333333
$ws.Range("K11").Value = 0.1958333333333333
$ws.Range("L11").Value = 0.5791666666666667
$ws.Range("S11").Value = 0.01666666666666667
$ws.Range("G12").Value = 0.7430555555555556
$ws.Range("J12").Value = 0.2152777777777778
$ws.Range("K12").Value = 0.006944444444444444
$ws.Range("L12").Value = 0.02083333333333333
$ws.Range("S12").Value = 0.01388888888888889
$ws.Range("G13").Value = 0.4791666666666667
$ws.Range("J13").Value = 0.5
$ws.Range("S13").Value = 0.02083333333333333
$ws.Range("F15").Value = 0.009345794392523364
$ws.Range("H15").Value = 0.1214953271028037
$ws.Range("I15").Value = 0.05607476635514019
$ws.Range("J15").Value = 0.4626168224299065
$ws.Range("K15").Value = 0.04205607476635514
$ws.Range("M15").Value = 0.009345794392523364
$ws.Range("O15").Value = 0.06074766355140187
$ws.Range("S15").Value = 0.2383177570093458
$ws.Range("F16").Value = 0.01036269430051814
$ws.Range("H16").Value = 0.1606217616580311
$ws.Range("I16").Value = 0.05699481865284974
$ws.Range("J16").Value = 0.4300518134715026
$ws.Range("K16").Value = 0.08808290155440414
$ws.Range("M16").Value = 0.03626943005181347
$ws.Range("O16").Value = 0.07253886010362694
$ws.Range("S16").Value = 0.1450777202072539
$ws.Range("F17").Value = 0.02558139534883721
$ws.Range("H17").Value = 0.1837209302325581
$ws.Range("I17").Value = 0.1023255813953488
$ws.Range("J17").Value = 0.4209302325581395
$ws.Range("K17").Value = 0.07209302325581396
$ws.Range("M17").Value = 0.01162790697674419
$ws.Range("N17").Value = 0.002325581395348837
$ws.Range("O17").Value = 0.07441860465116279
$ws.Range("S17").Value = 0.1069767441860465
$ws.Range("F18").Value = 0.02272727272727273
$ws.Range("H18").Value = 0.1875
$ws.Range("I18").Value = 0.07954545454545454
$ws.Range("J18").Value = 0.4147727272727273
$ws.Range("K18").Value = 0.1022727272727273
$ws.Range("M18").Value = 0.02272727272727273
$ws.Range("O18").Value = 0.09090909090909091
$ws.Range("S18").Value = 0.07954545454545454
$ws.Range("F19").Value = 0.01726973684210526
$ws.Range("H19").Value = 0.2105263157894737
$ws.Range("I19").Value = 0.09375
$ws.Range("J19").Value = 0.3963815789473684
$ws.Range("K19").Value = 0.09457236842105263
$ws.Range("M19").Value = 0.02467105263157895
$ws.Range("N19").Value = 0.001644736842105263
$ws.Range("O19").Value = 0.06990131578947369
$ws.Range("S19").Value = 0.0912828947368421
